$d = $word.ActiveDocument

function New-ParaXml([string]$innerPPr, [string]$innerRuns) {
    return @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$innerPPr$innerRuns</w:p>
</w:body>
</w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

$pPr0 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$pPr1 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$pPr2 = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# ---------------------------------------------------------------------------
# Step A: work on paragraph 7 (the bookmark paragraph) and everything that
# follows it, growing the list downward. Doing this before touching
# paragraphs 2/3/5 keeps those earlier indices stable.
# ---------------------------------------------------------------------------

# a1/a2: new paragraph after para 6 -> "Weiterhin ..." (ilvl=1)
$d.Paragraphs(6).Range.InsertParagraphAfter()
$runs = '<w:r><w:t>Weiterhin m' + [char]0xFC + 'ssen entsprechende Nachrichten an alle Benutzer (Session) eine Nachricht bekommen, dass ein Objekt gezeichnet wurde.</w:t></w:r>'
$d.Paragraphs(7).Range.InsertXML((New-ParaXml $pPr1 $runs))

# a3/a4: new paragraph after para 7 -> "Zudem ..." (ilvl=1)
$d.Paragraphs(7).Range.InsertParagraphAfter()
$runs = '<w:r><w:t>Zudem m' + [char]0xFC + 'ssen Nachrichten an alle gesendet werden, dass bspw. ein Objekt wieder gel' + [char]0xF6 + 'scht wurde.</w:t></w:r>'
$d.Paragraphs(8).Range.InsertXML((New-ParaXml $pPr1 $runs))

# a5/a6: new paragraph after para 8 -> "Hierdurch kam die Frage..." (ilvl=2)
$d.Paragraphs(8).Range.InsertParagraphAfter()
$runs = '<w:r><w:t>Hierdurch kam die Frage, wie die Implementierung sinnvoll aussieht, damit nur m' + [char]0xF6 + 'glichst wenig, aber effizienter Code produziert wird.</w:t></w:r>'
$d.Paragraphs(9).Range.InsertXML((New-ParaXml $pPr2 $runs))

# a7/a8: new paragraph after para 9 -> "Um dies zu bewerkstelligen..." (ilvl=2)
$d.Paragraphs(9).Range.InsertParagraphAfter()
$runs = '<w:r><w:t>Um dies zu bewerkstelligen, werden einfach immer dieselben Methoden (m' + [char]0xF6 + 'glichst vereinheitlicht) aufgerufen und die Buttons selbst rufen eine weitere Hilfs-Funktion auf, wodurch besagte Einheitsfunktion aufgerufen wird.</w:t></w:r>'
$d.Paragraphs(10).Range.InsertXML((New-ParaXml $pPr2 $runs))

# a9/a10: new paragraph after para 10 -> "Zuletzt sind Schwierigkeiten..." (ilvl=0)
$d.Paragraphs(10).Range.InsertParagraphAfter()
$runs = '<w:r><w:t>Zuletzt sind Schwierigkeiten aufgetreten bei der Speicherung bzw. dem Downloaden des entstandenen Bildes sowie der Historie.</w:t></w:r>'
$d.Paragraphs(11).Range.InsertXML((New-ParaXml $pPr0 $runs))

# a11: the bookmark paragraph is now paragraph 12. Edit it in place so the
# bookmark stays put: prepend a run of text, then append a "." run after
# the bookmark by inserting a helper paragraph and merging it back in.
$bookmarkPara = $d.Paragraphs(12)
$bookmarkPara.Range.InsertBefore('Hierzu gab es einige Hilfen im Internet, welche auf beide Szenarien angewendet werden konnten')

$bookmarkPara = $d.Paragraphs(12)
$endPoint = $d.Range($bookmarkPara.Range.End, $bookmarkPara.Range.End)
$runs = '<w:r><w:t>.</w:t></w:r>'
$endPoint.InsertXML((New-ParaXml $pPr1 $runs))

$bookmarkPara = $d.Paragraphs(12)
$mergeRange = $d.Range($bookmarkPara.Range.End - 1, $bookmarkPara.Range.End)
$mergeRange.Delete()

# ---------------------------------------------------------------------------
# Step B: paragraphs 2, 3 and 5 - merge runs and drop proofErr wrapping.
# ---------------------------------------------------------------------------

# Paragraph 2: "Grundlegende " + "Websocketimplementierung" -> single run,
# plus the new sentence about the tutorial with the Wingdings arrow symbol.
$runs = '<w:r><w:t>Grundlegende Websocketimplementierung</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> (Tutorial angeschaut und Beispiel nachprogrammiert </w:t></w:r>' +
        '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
        '<w:r><w:t xml:space="preserve"> danach nach eigenen Vorstellungen und W' + [char]0xFC + 'nschen umprogrammiert)</w:t></w:r>'
$d.Paragraphs(2).Range.InsertXML((New-ParaXml $pPr0 $runs))

# Paragraph 3: drop the proofErr wrap around "Nachrichtenübetragung".
$runs = '<w:r><w:t>Nachrichten' + [char]0xFC + 'betragung</w:t></w:r>'
$d.Paragraphs(3).Range.InsertXML((New-ParaXml $pPr0 $runs))

# Paragraph 5: merge the three runs ("Wie müssen..." + "JSONObjekt" + " aus?)?")
# into a single run, dropping the proofErr wrap.
$runs = '<w:r><w:t>Wie m' + [char]0xFC + 'ssen diese aussehen/aufgebaut sein (wie sieht das JSONObjekt aus?)?</w:t></w:r>'
$d.Paragraphs(5).Range.InsertXML((New-ParaXml $pPr1 $runs))

Write-Output "done"
